$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: I1 "AbilityText" -> "FavourText" (keep the quote-prefixed style) ---
$ws.Range("I1").Value = "'FavourText"

# --- Row 2: King ---
$ws.Range("A2").Value = "King"
$ws.Range("B2").Value = "King.ai"
$ws.Range("C2").Value = "Heart.ai"
$ws.Range("D2").Value = "Heart.ai"
$ws.Range("E2").Value = "Heart.ai"
$ws.Range("F2").Value = "Heart.ai"
$ws.Range("I2").Value = "Once per game:`n- Peek at a players role OR`n- Move a Delegate between Courts"

# --- Row 3: Traitor ---
$ws.Range("A3").Value = "Traitor"
$ws.Range("B3").Value = "Traitor.ai"
$ws.Range("D3").Value = "Heart.ai"
$ws.Range("E3").Value = "Heart.ai"
$ws.Range("G3").Value = "Damage.ai"
$ws.Range("H3").Value = "x2 to king, if alive."
$ws.Range("I3").Value = "Once per game:`n- Peek at a players role OR`n- Move a Delegate between Courts"

# --- Row 4: Loyalist ---
$ws.Range("A4").Value = "Loyalist"
$ws.Range("B4").Value = "Loyalist.ai"
$ws.Range("D4").Value = "Heart.ai"
$ws.Range("E4").Value = "Heart.ai"
$ws.Range("G4").Value = "Protection.ai"
$ws.Range("H4").Value = "x2 to king, if alive."
$ws.Range("I4").Value = "Once per game:`n- Peek at a players role OR`n- Move a Delegate between Courts"

# --- Formatting: build the needed styles in an order that doesn't leave orphaned xfs ---
# Style A: horizontal+vertical center (role-name column, A2:A4)
$ws.Range("A2:A4").HorizontalAlignment = -4108
$ws.Range("A2:A4").VerticalAlignment = -4108

# Style B: vertical center only (B2:H4, reuses the "vertical center" xf created above)
$ws.Range("B2:H4").VerticalAlignment = -4108

# Style C: vertical center + wrap (I2:I4)
$ws.Range("I2:I4").VerticalAlignment = -4108
$ws.Range("I2:I4").WrapText = $true

# Wrapping text auto-grew the row height; restore the original custom height.
$ws.Rows("2:4").RowHeight = 21.75

# --- Drop the now-unused trailing columns J:M (width range shrinks from 7-13 to 7-9) ---
$ws.Columns("J:M").Delete()

# --- Drop the now-unused trailing empty rows 5:20 ---
$ws.Rows("5:20").Delete()

# --- Match the saved selection from the source edit ---
$ws.Range("D14").Select()
